# Handback status report refresh: update the "Correspond Handoff Datetime"
# and "Correspond Handback DateTime" cells for the first file row on each
# language sheet (zh-cn and de-de), simulating a fresh report generation.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: row 2 corresponds to the
# 83f153c1-74e8-465e-9c49-b796e5acc545 file.
$wsZhCn.Range("E2").Value = "2016-03-22 04:46:07"
$wsZhCn.Range("H2").Value = "2016-03-22 04:46:29"

# de-de sheet: row 2 corresponds to the
# 83f153c1-74e8-465e-9c49-b796e5acc545 file.
$wsDeDe.Range("E2").Value = "2016-03-22 04:46:11"
$wsDeDe.Range("H2").Value = "2016-03-22 04:46:34"
